$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.615.94"
$ws.Range("E2").Value = "  -2.59%  "
$ws.Range("D3").Value = "2.231.02"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "297.72"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.48%  "
$ws.Range("E7").Value = "  -1.88%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "45.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -7.43%  "
$ws.Range("E11").Value = "  -2.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("E14").Value = "  -3.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.913"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("D17").Value = "2.566.95"
$ws.Range("E17").Value = "  -2.46%  "
$ws.Range("D18").Value = "2.251.33"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("D19").Value = "42.471.19"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.36"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.09%  "
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +23.37%  "
$ws.Range("E24").Value = "  -5.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "230.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.42%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  -1.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0888"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.78%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.72"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.77%  "
$ws.Range("E38").Value = "  -2.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0367"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("E40").Value = "  -2.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.52"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.80%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.96%  "
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.21%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.76"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -2.08%  "
